$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Feature'
$ws.Cells.Item(1, 2).Value = 'Importance'

$ws.Cells.Item(2, 1).Value = 'indeferiu'
$ws.Cells.Item(2, 2).Value = 0.0625
$ws.Cells.Item(3, 1).Value = 'deferimento'
$ws.Cells.Item(3, 2).Value = 0.0625
$ws.Cells.Item(4, 1).Value = 'concessão ordem'
$ws.Cells.Item(4, 2).Value = 0.046875
$ws.Cells.Item(5, 1).Value = 'concessão'
$ws.Cells.Item(5, 2).Value = 0.046875
$ws.Cells.Item(6, 1).Value = 'estado paulo'
$ws.Cells.Item(6, 2).Value = 0.03125
$ws.Cells.Item(7, 1).Value = 'óbice'
$ws.Cells.Item(7, 2).Value = 0.03125
$ws.Cells.Item(8, 1).Value = 'ordem prisão'
$ws.Cells.Item(8, 2).Value = 0.03125
$ws.Cells.Item(9, 1).Value = 'deferida'
$ws.Cells.Item(9, 2).Value = 0.03125
$ws.Cells.Item(10, 1).Value = 'inconstitucional'
$ws.Cells.Item(10, 2).Value = 0.03125
$ws.Cells.Item(11, 1).Value = 'sobrestamento'
$ws.Cells.Item(11, 2).Value = 0.03125
$ws.Cells.Item(12, 1).Value = 'senha relatório'
$ws.Cells.Item(12, 2).Value = 0.03125
$ws.Cells.Item(13, 1).Value = 'liminar hc'
$ws.Cells.Item(13, 2).Value = 0.03125
$ws.Cells.Item(14, 1).Value = 'deferi'
$ws.Cells.Item(14, 2).Value = 0.03125
$ws.Cells.Item(15, 1).Value = 'relatório http'
$ws.Cells.Item(15, 2).Value = 0.03125
$ws.Cells.Item(16, 1).Value = 'manifesta ilegalidade'
$ws.Cells.Item(16, 2).Value = 0.015625
$ws.Cells.Item(17, 1).Value = 'liminar suspender'
$ws.Cells.Item(17, 2).Value = 0.015625
$ws.Cells.Item(18, 1).Value = 'liminarmente pedido'
$ws.Cells.Item(18, 2).Value = 0.015625
$ws.Cells.Item(19, 1).Value = 'liminar deferida'
$ws.Cells.Item(19, 2).Value = 0.015625
$ws.Cells.Item(20, 1).Value = 'anotou'
$ws.Cells.Item(20, 2).Value = 0.015625
$ws.Cells.Item(21, 1).Value = 'indeferiu liminarmente'
$ws.Cells.Item(21, 2).Value = 0.015625
$ws.Cells.Item(22, 1).Value = 'opina concessão'
$ws.Cells.Item(22, 2).Value = 0.015625
$ws.Cells.Item(23, 1).Value = 'anos dois'
$ws.Cells.Item(23, 2).Value = 0.015625
$ws.Cells.Item(24, 1).Value = 'opina deferimento'
$ws.Cells.Item(24, 2).Value = 0.015625
$ws.Cells.Item(25, 1).Value = 'opinou concessão'
$ws.Cells.Item(25, 2).Value = 0.015625
$ws.Cells.Item(26, 1).Value = 'outro motivo'
$ws.Cells.Item(26, 2).Value = 0.015625
$ws.Cells.Item(27, 1).Value = 'preventiva fundamentos'
$ws.Cells.Item(27, 2).Value = 0.015625
$ws.Cells.Item(28, 1).Value = 'processo revelador'
$ws.Cells.Item(28, 2).Value = 0.015625
$ws.Cells.Item(29, 1).Value = 'providência procuradoria'
$ws.Cells.Item(29, 2).Value = 0.015625
$ws.Cells.Item(30, 1).Value = 'república concessão'
$ws.Cells.Item(30, 2).Value = 0.015625
$ws.Cells.Item(31, 1).Value = 'senha'
$ws.Cells.Item(31, 2).Value = 0.015625
$ws.Cells.Item(32, 1).Value = 'stj indeferiu'
$ws.Cells.Item(32, 2).Value = 0.015625
$ws.Cells.Item(33, 1).Value = 'substituir prisão'
$ws.Cells.Item(33, 2).Value = 0.015625
$ws.Cells.Item(34, 1).Value = 'violência grave'
$ws.Cells.Item(34, 2).Value = 0.015625
$ws.Cells.Item(35, 1).Value = 'xliii constituição'
$ws.Cells.Item(35, 2).Value = 0.015625
$ws.Cells.Item(36, 1).Value = 'presa'
$ws.Cells.Item(36, 2).Value = 0.015625
$ws.Cells.Item(37, 1).Value = 'justiça nº'
$ws.Cells.Item(37, 2).Value = 0.015625
$ws.Cells.Item(38, 1).Value = 'contornos'
$ws.Cells.Item(38, 2).Value = 0.015625
$ws.Cells.Item(39, 1).Value = 'cautelares previstas'
$ws.Cells.Item(39, 2).Value = 0.015625
$ws.Cells.Item(40, 1).Value = 'deste habeas'
$ws.Cells.Item(40, 2).Value = 0.015625
$ws.Cells.Item(41, 1).Value = 'deferimento liminar'
$ws.Cells.Item(41, 2).Value = 0.015625
$ws.Cells.Item(42, 1).Value = 'causas aumento'
$ws.Cells.Item(42, 2).Value = 0.015625
$ws.Cells.Item(43, 1).Value = 'ficou'
$ws.Cells.Item(43, 2).Value = 0.015625
$ws.Cells.Item(44, 1).Value = 'enunciado súmula'
$ws.Cells.Item(44, 2).Value = 0.015625
$ws.Cells.Item(45, 1).Value = 'campo precário'
$ws.Cells.Item(45, 2).Value = 0.015625
$ws.Cells.Item(46, 1).Value = 'publiquem'
$ws.Cells.Item(46, 2).Value = 0.0
$ws.Cells.Item(47, 1).Value = 'colham'
$ws.Cells.Item(47, 2).Value = 0.0
$ws.Cells.Item(48, 1).Value = 'senha primeira'
$ws.Cells.Item(48, 2).Value = 0.0
$ws.Cells.Item(49, 1).Value = 'república parecer'
$ws.Cells.Item(49, 2).Value = 0.0
$ws.Cells.Item(50, 1).Value = 'resumida'
$ws.Cells.Item(50, 2).Value = 0.0
$ws.Cells.Item(51, 1).Value = 'resumida prisão'
$ws.Cells.Item(51, 2).Value = 0.0
$ws.Cells.Item(52, 1).Value = 'revelador'
$ws.Cells.Item(52, 2).Value = 0.0
$ws.Cells.Item(53, 1).Value = 'revelou contornos'
$ws.Cells.Item(53, 2).Value = 0.0
$ws.Cells.Item(54, 1).Value = 'proferida ministro'
$ws.Cells.Item(54, 2).Value = 0.0
$ws.Cells.Item(55, 1).Value = 'procuradoria geral'
$ws.Cells.Item(55, 2).Value = 0.0
$ws.Cells.Item(56, 1).Value = 'ficou assim'
$ws.Cells.Item(56, 2).Value = 0.0
$ws.Cells.Item(57, 1).Value = 'procuradoria'
$ws.Cells.Item(57, 2).Value = 0.0
$ws.Cells.Item(58, 1).Value = 'ser julgado'
$ws.Cells.Item(58, 2).Value = 0.0
$ws.Cells.Item(59, 1).Value = 'sob código'
$ws.Cells.Item(59, 2).Value = 0.0
$ws.Cells.Item(60, 1).Value = 'aurélio decisão'
$ws.Cells.Item(60, 2).Value = 0.0
$ws.Cells.Item(61, 1).Value = 'assim revelou'
$ws.Cells.Item(61, 2).Value = 0.0
$ws.Cells.Item(62, 1).Value = 'substituir'
$ws.Cells.Item(62, 2).Value = 0.0
$ws.Cells.Item(63, 1).Value = 'assim resumida'
$ws.Cells.Item(63, 2).Value = 0.0
$ws.Cells.Item(64, 1).Value = 'suspender'
$ws.Cells.Item(64, 2).Value = 0.0
$ws.Cells.Item(65, 1).Value = 'suspender efeitos'
$ws.Cells.Item(65, 2).Value = 0.0
$ws.Cells.Item(66, 1).Value = 'teixeira'
$ws.Cells.Item(66, 2).Value = 0.0
$ws.Cells.Item(67, 1).Value = 'vedação liberdade'
$ws.Cells.Item(67, 2).Value = 0.0
$ws.Cells.Item(68, 1).Value = 'assessoria prestou'
$ws.Cells.Item(68, 2).Value = 0.0
$ws.Cells.Item(69, 1).Value = 'assessoria'
$ws.Cells.Item(69, 2).Value = 0.0
$ws.Cells.Item(70, 1).Value = 'análise pedido'
$ws.Cells.Item(70, 2).Value = 0.0
$ws.Cells.Item(71, 1).Value = 'brasília residência'
$ws.Cells.Item(71, 2).Value = 0.0
$ws.Cells.Item(72, 1).Value = 'prestadas gabinete'
$ws.Cells.Item(72, 2).Value = 0.0
$ws.Cells.Item(73, 1).Value = 'colham parecer'
$ws.Cells.Item(73, 2).Value = 0.0
$ws.Cells.Item(74, 1).Value = 'liminar assessor'
$ws.Cells.Item(74, 2).Value = 0.0
$ws.Cells.Item(75, 1).Value = 'gabinete prestou'
$ws.Cells.Item(75, 2).Value = 0.0
$ws.Cells.Item(76, 1).Value = 'habeas contra'
$ws.Cells.Item(76, 2).Value = 0.0
$ws.Cells.Item(77, 1).Value = 'idêntica medida'
$ws.Cells.Item(77, 2).Value = 0.0
$ws.Cells.Item(78, 1).Value = 'impetração eis'
$ws.Cells.Item(78, 2).Value = 0.0
$ws.Cells.Item(79, 1).Value = 'implicou deferimento'
$ws.Cells.Item(79, 2).Value = 0.0
$ws.Cells.Item(80, 1).Value = 'espécie ficou'
$ws.Cells.Item(80, 2).Value = 0.0
$ws.Cells.Item(81, 1).Value = 'informado'
$ws.Cells.Item(81, 2).Value = 0.0
$ws.Cells.Item(82, 1).Value = 'informado análise'
$ws.Cells.Item(82, 2).Value = 0.0
$ws.Cells.Item(83, 1).Value = 'informações paciente'
$ws.Cells.Item(83, 2).Value = 0.0
$ws.Cells.Item(84, 1).Value = 'julgado turma'
$ws.Cells.Item(84, 2).Value = 0.0
$ws.Cells.Item(85, 1).Value = 'enunciado'
$ws.Cells.Item(85, 2).Value = 0.0
$ws.Cells.Item(86, 1).Value = 'juízo criminal'
$ws.Cells.Item(86, 2).Value = 0.0
$ws.Cells.Item(87, 1).Value = 'eis informado'
$ws.Cells.Item(87, 2).Value = 0.0
$ws.Cells.Item(88, 1).Value = 'fundamentos insubsistência'
$ws.Cells.Item(88, 2).Value = 0.0
$ws.Cells.Item(89, 1).Value = 'liminar espécie'
$ws.Cells.Item(89, 2).Value = 0.0
$ws.Cells.Item(90, 1).Value = 'efêmero'
$ws.Cells.Item(90, 2).Value = 0.0
$ws.Cells.Item(91, 1).Value = 'deferida assessoria'
$ws.Cells.Item(91, 2).Value = 0.0
$ws.Cells.Item(92, 1).Value = 'deferi pedido'
$ws.Cells.Item(92, 2).Value = 0.0
$ws.Cells.Item(93, 1).Value = 'decretada desfavor'
$ws.Cells.Item(93, 2).Value = 0.0
$ws.Cells.Item(94, 1).Value = 'decisão proferida'
$ws.Cells.Item(94, 2).Value = 0.0
$ws.Cells.Item(95, 1).Value = 'decisão implicou'
$ws.Cells.Item(95, 2).Value = 0.0
$ws.Cells.Item(96, 1).Value = 'código senha'
$ws.Cells.Item(96, 2).Value = 0.0
$ws.Cells.Item(97, 1).Value = 'paulo indeferiu'
$ws.Cells.Item(97, 2).Value = 0.0
$ws.Cells.Item(98, 1).Value = 'precário'
$ws.Cells.Item(98, 2).Value = 0.0
$ws.Cells.Item(99, 1).Value = 'precário efêmero'
$ws.Cells.Item(99, 2).Value = 0.0
$ws.Cells.Item(100, 1).Value = 'contornos impetração'
$ws.Cells.Item(100, 2).Value = 0.0
$ws.Cells.Item(101, 1).Value = 'óbice previsto'
$ws.Cells.Item(101, 2).Value = 0.0

Write-Host "done"
